# Insert a new price record for "Vega Monumental Concepción" (Ajo / Chino / Primera)
# as row 180, pushing the existing rows 180-291 down to 181-292.
#
# This mirrors the weekly update pattern seen in the commit: a brand-new
# observation is prepended near the top of the data block and every older
# row shifts down by one, with the previously-last row (old 291) ending up
# at the new last row (292).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 180..291 down by one, leaving a blank row 180 to fill in.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new observation.
$ws.Cells.Item(180, 1).Value  = 11
$ws.Cells.Item(180, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(180, 3).Value  = "Bíobío"
$ws.Cells.Item(180, 4).Value  = 45086
$ws.Cells.Item(180, 5).Value  = 8
$ws.Cells.Item(180, 6).Value  = 100112003
$ws.Cells.Item(180, 7).Value  = "Ajo"
$ws.Cells.Item(180, 8).Value  = "Chino"
$ws.Cells.Item(180, 9).Value  = "Primera"
$ws.Cells.Item(180, 10).Value = 150
$ws.Cells.Item(180, 11).Value = 15000
$ws.Cells.Item(180, 12).Value = 15000
$ws.Cells.Item(180, 13).Value = 15000
$ws.Cells.Item(180, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(180, 15).Value = "China"
$ws.Cells.Item(180, 16).Value = 1500
$ws.Cells.Item(180, 17).Value = 10
$ws.Cells.Item(180, 18).Value = "Hortaliza"
